# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" conversion summary text ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.67 = 6261.02 pesos`n✅ 6261.02 pesos = 1.65 = 906.92 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet rate values ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("O10").Value = 3749.35
$tasas.Range("N12").Value = 3796.99
$tasas.Range("O12").Value = 550
